$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DPbES")

# Update the "environmentally preferred" example dispatch-priority values
# (column B) for the rows whose numbers changed in this revision. Every
# other cell in the row is a formula copied across (shared formula) from
# column B, so the engine's recalculation will refresh C:AK automatically.

# hard coal: 5 -> 2
$ws.Range("B2").Value = 2

# natural gas nonpeaker: 4 -> 1
$ws.Range("B3").Value = 1

# hydro: 2 -> 5
$ws.Range("B5").Value = 5

# biomass: 3 -> 5
$ws.Range("B9").Value = 5

# crude oil: was "=B11" (1) -> literal 5 (formula link removed)
$ws.Range("B15").Value = 5

# heavy or residual fuel oil: was "=B11" (1) -> literal 5 (formula link removed)
$ws.Range("B16").Value = 5

# municipal solid waste: was "=B9" (3) -> literal 3, and the rest of the row
# (C17:AK17), which used to mirror row 9, now references B17 directly.
$ws.Range("B17").Value = 3
$ws.Range("C17").Formula = "=`$B`$17"
$ws.Range("D17:AK17").Formula = "=`$B`$17"

# Reflect the saved UI state: the DPbES sheet is the active tab, with I17
# selected (the About sheet loses its "tabSelected" flag as a result).
$ws.Activate()
$ws.Range("I17").Select()
